$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.029.64"
$ws.Range("E2").Value = "  -0.61%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.952.14"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.52%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.46"
$ws.Range("E5").Value = "  -2.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.43"
$ws.Range("E6").Value = "  +1.42%  "

$ws.Range("E7").Value = "  +0.38%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.943.45"
$ws.Range("E8").Value = "  -1.29%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -2.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  +10.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.144"
$ws.Range("E11").Value = "  -2.32%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("E12").Value = "  -0.75%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000223"
$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.19"
$ws.Range("E14").Value = "  +0.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.125"
$ws.Range("E15").Value = "  -1.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.469.62"
$ws.Range("E16").Value = "  -0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.416.02"
$ws.Range("E17").Value = "  +0.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.84"
$ws.Range("E18").Value = "  -0.97%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.973.10"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.91"
$ws.Range("E20").Value = "  -2.64%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.78"
$ws.Range("E21").Value = "  -0.79%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.673"
$ws.Range("E22").Value = "  -1.09%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.30"
$ws.Range("E23").Value = "  -0.21%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.61"
$ws.Range("E24").Value = "  -0.54%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.04"
$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  +0.26%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.86"
$ws.Range("E27").Value = "  -1.12%  "

$ws.Range("E28").Value = "  -0.07%  "

$ws.Range("B29").Value = "FirstDigitalUSD"
$ws.Range("C29").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.33%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.66"
$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.24"
$ws.Range("E31").Value = "  +0.09%  "

$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.19"
$ws.Range("E32").Value = "  +6.96%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.78"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("E34").Value = "  -2.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0822"
$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  -0.49%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.68"
$ws.Range("E37").Value = "  -1.20%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.92"
$ws.Range("E38").Value = "  -0.53%  "

$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.96"
$ws.Range("E39").Value = "  +4.10%  "

$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.125"
$ws.Range("E40").Value = "  +1.17%  "

$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.02"
$ws.Range("E41").Value = "  +0.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").Value = "  -2.20%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.294"
$ws.Range("E43").Value = "  +9.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.70"
$ws.Range("E44").Value = "  +5.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0348"
$ws.Range("E45").Value = "  -0.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "369.60"
$ws.Range("E46").Value = "  -4.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.654.06"
$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.42"
$ws.Range("E48").Value = "  -0.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.51"
$ws.Range("E49").Value = "  +9.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("E51").Value = "  +1.07%  "

